$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 15421.25
$ws.Range("J13").Value = 15595
$ws.Range("L13").Value = 15595
$ws.Range("N13").Value = -15933
$ws.Range("H17").Value = 2291769.8
$ws.Range("J17").Value = 2353539.2
$ws.Range("L17").Value = 7060617.600000001
$ws.Range("N17").Value = -7060953.600000001
$ws.Range("H48").Value = 2008.5
$ws.Range("I48").Value = 2008.5
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 6025.5
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -5733.5
$ws.Range("N48").ClearContents()
$ws.Range("H56").Value = 2008.5
$ws.Range("I56").Value = 2008.5
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 6025.5
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -5491.5
$ws.Range("N56").ClearContents()
$ws.Range("H112").Value = 2343.9524
$ws.Range("J112").Value = 2486.6316
$ws.Range("L112").Value = 7459.8948
$ws.Range("N112").Value = -9675.8948

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 2496.5
$ws.Range("J17").Value = 2495
$ws.Range("L17").Value = 2495
$ws.Range("N17").Value = -2841
$ws.Range("H22").Value = 3289.375
$ws.Range("I22").Value = 2463.2
$ws.Range("J22").Value = 4666.3335
$ws.Range("K22").Value = 2463.2
$ws.Range("L22").Value = 4666.3335
$ws.Range("M22").Value = -2164.2
$ws.Range("N22").Value = -5264.3335
$ws.Range("H27").Value = 27299.666
$ws.Range("J27").Value = 27299.666
$ws.Range("L27").Value = 27299.666
$ws.Range("N27").Value = -27667.666
$ws.Range("H35").Value = 15077.5
$ws.Range("I35").Value = 6157.2
$ws.Range("K35").Value = 6157.2
$ws.Range("M35").Value = -5751.2
$ws.Range("H39").Value = 15016
$ws.Range("J39").Value = 15016
$ws.Range("L39").Value = 15016
$ws.Range("N39").Value = -16056
$ws.Range("H40").Value = 74995
$ws.Range("J40").Value = 99993.5
$ws.Range("L40").Value = 99993.5
$ws.Range("N40").Value = -100345.5
$ws.Range("H49").Value = 65513
$ws.Range("J49").Value = 65513
$ws.Range("L49").Value = 65513
$ws.Range("N49").Value = -66033

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H32").Value = 25464.5
$ws.Range("J32").Value = 25464.5
$ws.Range("L32").Value = 25464.5
$ws.Range("N32").Value = -26232.5
$ws.Range("H36").Value = 10040
$ws.Range("I36").Value = 10040
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 10040
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -9506
$ws.Range("N36").ClearContents()
$ws.Range("H38").Value = 39899
$ws.Range("J38").Value = 39899
$ws.Range("L38").Value = 39899
$ws.Range("N38").Value = -40731

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 5185.4287
$ws.Range("I8").Value = 1999.5
$ws.Range("J8").Value = 6459.8
$ws.Range("K8").Value = 1999.5
$ws.Range("L8").Value = 6459.8
$ws.Range("M8").Value = -1859.5
$ws.Range("N8").Value = -6739.8
$ws.Range("H22").Value = 472.26086
$ws.Range("I22").Value = 219.1579
$ws.Range("K22").Value = 219.1579
$ws.Range("M22").Value = 130.8421
$ws.Range("H23").Value = 20684.111
$ws.Range("I23").Value = 19296.334
$ws.Range("J23").Value = 21378
$ws.Range("K23").Value = 19296.334
$ws.Range("L23").Value = 21378
$ws.Range("M23").Value = -19056.334
$ws.Range("N23").Value = -21858
$ws.Range("H27").Value = 20684.111
$ws.Range("I27").Value = 19296.334
$ws.Range("J27").Value = 21378
$ws.Range("K27").Value = 19296.334
$ws.Range("L27").Value = 21378
$ws.Range("M27").Value = -19104.334
$ws.Range("N27").Value = -21762
$ws.Range("H75").Value = 45997
$ws.Range("J75").Value = 45997
$ws.Range("L75").Value = 45997
$ws.Range("N75").Value = -47993
$ws.Range("H78").Value = 45997
$ws.Range("J78").Value = 45997
$ws.Range("L78").Value = 137991
$ws.Range("N78").Value = -147975

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 45.633335
$ws.Range("I2").Value = 50
$ws.Range("J2").Value = 33.625
$ws.Range("K2").Value = 300
$ws.Range("L2").Value = 201.75
$ws.Range("M2").Value = -187
$ws.Range("N2").Value = -427.75
$ws.Range("H109").Value = 4000
$ws.Range("I109").Value = 4000
$ws.Range("K109").Value = 12000
$ws.Range("M109").Value = -10960
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 549998
$ws.Range("I11").Value = 549998
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 549998
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -549859
$ws.Range("N11").ClearContents()
$ws.Range("H22").Value = 400
$ws.Range("J22").Value = 366.66666
$ws.Range("L22").Value = 366.66666
$ws.Range("N22").Value = -1424.66666
$ws.Range("H23").Value = 716
$ws.Range("J23").Value = 716
$ws.Range("L23").Value = 716
$ws.Range("N23").Value = -1162
$ws.Range("H25").Value = 1866.6666
$ws.Range("J25").Value = 1866.6666
$ws.Range("L25").Value = 1866.6666
$ws.Range("N25").Value = -2924.6666
$ws.Range("H33").Value = 39799
$ws.Range("J33").Value = 44699.5
$ws.Range("L33").Value = 44699.5
$ws.Range("N33").Value = -45203.5
$ws.Range("H36").Value = 34230
$ws.Range("I36").Value = 2690
$ws.Range("J36").Value = 50000
$ws.Range("K36").Value = 2690
$ws.Range("L36").Value = 50000
$ws.Range("M36").Value = -2205
$ws.Range("N36").Value = -50970

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 2009
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H19").Value = 5663.125
$ws.Range("J19").Value = 6866.5
$ws.Range("L19").Value = 6866.5
$ws.Range("N19").Value = -7206.5
$ws.Range("H20").Value = 95000.22
$ws.Range("J20").Value = 97272.95
$ws.Range("L20").Value = 97272.95
$ws.Range("N20").Value = -97724.95
$ws.Range("H31").Value = 7713.4287
$ws.Range("J31").Value = 10218.8
$ws.Range("L31").Value = 10218.8
$ws.Range("N31").Value = -10714.8
$ws.Range("H32").Value = 3335337.2
$ws.Range("I32").Value = 3335337.2
$ws.Range("K32").Value = 3335337.2
$ws.Range("M32").Value = -3335020.2
$ws.Range("H61").Value = 2649.889
$ws.Range("I61").Value = 2649.889
$ws.Range("K61").Value = 2649.889
$ws.Range("M61").Value = -2447.889
$ws.Range("H113").Value = 2649.889
$ws.Range("I113").Value = 2649.889
$ws.Range("K113").Value = 2649.889
$ws.Range("M113").Value = -479.8890000000001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 25666.666
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H26").Value = 18255.5
$ws.Range("I26").Value = 17670.666
$ws.Range("J26").Value = 20010
$ws.Range("K26").Value = 17670.666
$ws.Range("L26").Value = 20010
$ws.Range("M26").Value = -17377.666
$ws.Range("N26").Value = -20596
$ws.Range("H34").Value = 46914
$ws.Range("J34").Value = 55876.668
$ws.Range("L34").Value = 55876.668
$ws.Range("N34").Value = -56282.668
$ws.Range("H37").Value = 74991.5
$ws.Range("J37").Value = 99984
$ws.Range("L37").Value = 99984
$ws.Range("N37").Value = -100390
$ws.Range("H38").Value = 15000
$ws.Range("I38").Value = 15000
$ws.Range("K38").Value = 15000
$ws.Range("M38").Value = -14527
$ws.Range("H49").Value = 86248.75
$ws.Range("I49").Value = 44998.5
$ws.Range("J49").Value = 99998.836
$ws.Range("K49").Value = 44998.5
$ws.Range("L49").Value = 99998.836
$ws.Range("M49").Value = -44768.5
$ws.Range("N49").Value = -100458.836
$ws.Range("H105").Value = 101250
$ws.Range("J105").Value = 101250
$ws.Range("L105").Value = 101250
$ws.Range("N105").Value = -108238

Write-Host "Done updating Twintania_Profits sheets"